$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects")

# Rename "People Wearing Nehru Jacket " -> "Nehru Jacket " (naming convention cleanup)
# and add new category rows below the existing ones.
# Order below reproduces the exact shared-string insertion order of the
# original edit session.
$ws.Range("A5").Value = "Smiling Face"
$ws.Range("A2").Value = "Nehru Jacket "
$ws.Range("A4").Value = "Glasses"
$ws.Range("A6").Value = "Tree"
$ws.Range("A7").Value = "Birds"
$ws.Range("A8").Value = "Hat"
$ws.Range("A9").Value = "Bike"
$ws.Range("A10").Value = "Car"

$ws.Range("A11").Select()
